$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 11: MAgPIE 4.2 stand-alone AR6 registration form.
# Force column A to text first so the dd.mm.yyyy-style date string isn't
# auto-converted into a date serial number, then drop back to the
# workbook's default ("Normal") style so no stray number-format sticks to
# the cell.
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "08.07.2020"
$ws.Range("A11").Style = "Normal"

$ws.Range("B11").Value = "MAgPIE 4.2"
$ws.Range("C11").Value = "IPCC_AR6_model_registration_MAgPIE_4.2.xlsx"
$ws.Range("D11").Value = "Florian Humpenöder"

$ws.Range("E11").Value = "MAgPIE 4.2 stand-alone version used for peatland paper"
# Column E carries a wrap-text style by default; the new row's cell keeps
# the workbook default style instead, so reset it explicitly.
$ws.Range("E11").Style = "Normal"
